$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "DIAS" (days) values in column J for rows 2-10 (each increased by 4 days)
$ws.Range("J2").Value = 174
$ws.Range("J3").Value = 174
$ws.Range("J4").Value = 167
$ws.Range("J5").Value = 159
$ws.Range("J6").Value = 144
$ws.Range("J7").Value = 140
$ws.Range("J8").Value = 133
$ws.Range("J9").Value = 123
$ws.Range("J10").Value = 105
